$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")
$ws.Range("C2:C7").Value = "common"
$ws.Range("G3").Select()
